# Auto-generated edit script: adds a "2022-Q1" sheet with fund holdings data
# (mirroring the layout of the existing "2021-Q4" sheet), and records the new
# quarter as a summary row in the "总计" (grand-total) sheet.
#
# Strategy: write all cell *values* first (forcing text typing on numeric-
# looking strings via a temporary "@" number format so Excel does not
# silently coerce "26.19" / "004666" into numbers and lose precision/zeros),
# then in a second pass paste-special just the *formatting* from known-good
# template cells on top - this restores the exact look (bold/border header,
# bold/border index column, completely unstyled data cells) without the
# "@" format leaving a stray style behind.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right after "2021-Q4" and
#    before "总计" (target sheet order: 2021-Q4, 2022-Q1, 总计)
# ------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4Sheet)
$newSheet.Name = "2022-Q1"

# ---- header row (row 1) ----
Set-TextCell $newSheet.Range("B1") "基金代码"
Set-TextCell $newSheet.Range("C1") "基金名称"
Set-TextCell $newSheet.Range("D1") "基金规模"
Set-TextCell $newSheet.Range("E1") "股票总仓位"
Set-TextCell $newSheet.Range("F1") "仓位占比"
Set-TextCell $newSheet.Range("G1") "持有市值(亿元)"
Set-TextCell $newSheet.Range("H1") "仓位排名"

# ---- data rows (row 2 .. row 13) ----
$newSheet.Cells.Item(2,1).Value = 0
Set-TextCell $newSheet.Cells.Item(2,2) "004666"
Set-TextCell $newSheet.Cells.Item(2,3) "长城久嘉创新成长灵活配置混合"
Set-TextCell $newSheet.Cells.Item(2,4) "26.19"
Set-TextCell $newSheet.Cells.Item(2,5) "92.18"
Set-TextCell $newSheet.Cells.Item(2,6) "2.88"
Set-TextCell $newSheet.Cells.Item(2,7) "0.7543"
$newSheet.Cells.Item(2,8).Value = 8

$newSheet.Cells.Item(3,1).Value = 1
Set-TextCell $newSheet.Cells.Item(3,2) "001150"
Set-TextCell $newSheet.Cells.Item(3,3) "融通互联网传媒灵活配置混合"
Set-TextCell $newSheet.Cells.Item(3,4) "9.18"
Set-TextCell $newSheet.Cells.Item(3,5) "92.66"
Set-TextCell $newSheet.Cells.Item(3,6) "2.68"
Set-TextCell $newSheet.Cells.Item(3,7) "0.2460"
$newSheet.Cells.Item(3,8).Value = 8

$newSheet.Cells.Item(4,1).Value = 2
Set-TextCell $newSheet.Cells.Item(4,2) "013067"
Set-TextCell $newSheet.Cells.Item(4,3) "富安达中小盘六个月持有期混合"
Set-TextCell $newSheet.Cells.Item(4,4) "2.45"
Set-TextCell $newSheet.Cells.Item(4,5) "74.39"
Set-TextCell $newSheet.Cells.Item(4,6) "3.49"
Set-TextCell $newSheet.Cells.Item(4,7) "0.0855"
$newSheet.Cells.Item(4,8).Value = 4

$newSheet.Cells.Item(5,1).Value = 3
Set-TextCell $newSheet.Cells.Item(5,2) "004314"
Set-TextCell $newSheet.Cells.Item(5,3) "前海开源沪港深新硬件主题灵活配置混合A"
Set-TextCell $newSheet.Cells.Item(5,4) "1.67"
Set-TextCell $newSheet.Cells.Item(5,5) "90.05"
Set-TextCell $newSheet.Cells.Item(5,6) "3.65"
Set-TextCell $newSheet.Cells.Item(5,7) "0.0610"
$newSheet.Cells.Item(5,8).Value = 10

$newSheet.Cells.Item(6,1).Value = 4
Set-TextCell $newSheet.Cells.Item(6,2) "519113"
Set-TextCell $newSheet.Cells.Item(6,3) "浦银安盛精致生活混合"
Set-TextCell $newSheet.Cells.Item(6,4) "2.09"
Set-TextCell $newSheet.Cells.Item(6,5) "90.20"
Set-TextCell $newSheet.Cells.Item(6,6) "2.31"
Set-TextCell $newSheet.Cells.Item(6,7) "0.0483"
$newSheet.Cells.Item(6,8).Value = 10

$newSheet.Cells.Item(7,1).Value = 5
Set-TextCell $newSheet.Cells.Item(7,2) "004315"
Set-TextCell $newSheet.Cells.Item(7,3) "前海开源沪港深新硬件主题灵活配置混合C"
Set-TextCell $newSheet.Cells.Item(7,4) "1.00"
Set-TextCell $newSheet.Cells.Item(7,5) "90.05"
Set-TextCell $newSheet.Cells.Item(7,6) "3.65"
Set-TextCell $newSheet.Cells.Item(7,7) "0.0365"
$newSheet.Cells.Item(7,8).Value = 10

$newSheet.Cells.Item(8,1).Value = 6
Set-TextCell $newSheet.Cells.Item(8,2) "006692"
Set-TextCell $newSheet.Cells.Item(8,3) "金信消费升级股票A"
Set-TextCell $newSheet.Cells.Item(8,4) "0.62"
Set-TextCell $newSheet.Cells.Item(8,5) "94.05"
Set-TextCell $newSheet.Cells.Item(8,6) "4.59"
Set-TextCell $newSheet.Cells.Item(8,7) "0.0285"
$newSheet.Cells.Item(8,8).Value = 4

$newSheet.Cells.Item(9,1).Value = 7
Set-TextCell $newSheet.Cells.Item(9,2) "002292"
Set-TextCell $newSheet.Cells.Item(9,3) "诺安益鑫灵活配置混合"
Set-TextCell $newSheet.Cells.Item(9,4) "0.30"
Set-TextCell $newSheet.Cells.Item(9,5) "50.08"
Set-TextCell $newSheet.Cells.Item(9,6) "4.58"
Set-TextCell $newSheet.Cells.Item(9,7) "0.0137"
$newSheet.Cells.Item(9,8).Value = 3

$newSheet.Cells.Item(10,1).Value = 8
Set-TextCell $newSheet.Cells.Item(10,2) "002862"
Set-TextCell $newSheet.Cells.Item(10,3) "金信量化精选灵活配置混合"
Set-TextCell $newSheet.Cells.Item(10,4) "0.16"
Set-TextCell $newSheet.Cells.Item(10,5) "94.28"
Set-TextCell $newSheet.Cells.Item(10,6) "7.03"
Set-TextCell $newSheet.Cells.Item(10,7) "0.0112"
$newSheet.Cells.Item(10,8).Value = 1

$newSheet.Cells.Item(11,1).Value = 9
Set-TextCell $newSheet.Cells.Item(11,2) "005437"
Set-TextCell $newSheet.Cells.Item(11,3) "易方达易百智能量化策略灵活配置混合A"
Set-TextCell $newSheet.Cells.Item(11,4) "0.96"
Set-TextCell $newSheet.Cells.Item(11,5) "94.52"
Set-TextCell $newSheet.Cells.Item(11,6) "1.13"
Set-TextCell $newSheet.Cells.Item(11,7) "0.0108"
$newSheet.Cells.Item(11,8).Value = 4

$newSheet.Cells.Item(12,1).Value = 10
Set-TextCell $newSheet.Cells.Item(12,2) "006693"
Set-TextCell $newSheet.Cells.Item(12,3) "金信消费升级股票C"
Set-TextCell $newSheet.Cells.Item(12,4) "0.20"
Set-TextCell $newSheet.Cells.Item(12,5) "94.05"
Set-TextCell $newSheet.Cells.Item(12,6) "4.59"
Set-TextCell $newSheet.Cells.Item(12,7) "0.0092"
$newSheet.Cells.Item(12,8).Value = 4

$newSheet.Cells.Item(13,1).Value = 11
Set-TextCell $newSheet.Cells.Item(13,2) "005438"
Set-TextCell $newSheet.Cells.Item(13,3) "易方达易百智能量化策略灵活配置混合C"
Set-TextCell $newSheet.Cells.Item(13,4) "0.17"
Set-TextCell $newSheet.Cells.Item(13,5) "94.52"
Set-TextCell $newSheet.Cells.Item(13,6) "1.13"
Set-TextCell $newSheet.Cells.Item(13,7) "0.0019"
$newSheet.Cells.Item(13,8).Value = 4

# ---- restore formatting on the new sheet (paste-special, formats only) ----
# Header row + index column: bold / centered / thin-bordered, copied from "2021-Q4".
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2:A13").PasteSpecial(-4122) | Out-Null
# Plain data cells: no special style at all (also clears the temporary "@" format).
$q4Sheet.Range("B2").Copy()
$newSheet.Range("B2:G13").PasteSpecial(-4122) | Out-Null
# Rank column (H): plain/unstyled, like in "2021-Q4".
$q4Sheet.Range("H2").Copy()
$newSheet.Range("H2:H13").PasteSpecial(-4122) | Out-Null
$newSheet.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# 2. Update the "总计" (grand-total) sheet: insert a new row for 2022-Q1
#    above the existing 2021-Q4 row, shifting it down.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2,1).Value = 0
Set-TextCell $totalSheet.Cells.Item(2,2) "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 12
$totalSheet.Cells.Item(2,4).Value = 1.31

# Restore formatting on the new row: index-style column A (bold/bordered,
# copied from row 3), plain/unstyled B:D (clears the "@" format + any
# formatting the row-insert copied down from the header row).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122) | Out-Null
$totalSheet.Range("B3").Copy()
$totalSheet.Range("B2").PasteSpecial(-4122) | Out-Null
$totalSheet.Range("C3:D3").Copy()
$totalSheet.Range("C2:D2").PasteSpecial(-4122) | Out-Null
$totalSheet.Application.CutCopyMode = $false

# Keep the (now shifted) 2021-Q4 rows index value consistent (A3 = 1)
$totalSheet.Cells.Item(3,1).Value = 1

Write-Host "Done: added 2022-Q1 sheet and updated total-summary sheet."
